{"js": "// Update the date line and the 25 division problems in the practice\n// table. Replacements are applied by the paragraph's position in the\n// document (title paragraph + table cell paragraphs, in reading\n// order) rather than by text search, because several of the original\n// problem strings repeat (e.g. \"86\u00f73=\", \"41\u00f79=\", \"24\u00f78=\") and a naive\n// text replace would not let us target each occurrence independently.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// index -> [expected old text, new text]; index is the paragraph's\n// position within context.document.body.paragraphs.\nconst replacements = {\n  0: [\"2024-01-23 Tuesday\", \"2024-01-24 Wednesday\"],\n  1: [\"93\u00f77=\", \"40\u00f79=\"],\n  2: [\"86\u00f73=\", \"36\u00f77=\"],\n  3: [\"43\u00f79=\", \"19\u00f75=\"],\n  4: [\"98\u00f73=\", \"21\u00f78=\"],\n  5: [\"78\u00f73=\", \"82\u00f72=\"],\n  21: [\"78\u00f78=\", \"74\u00f74=\"],\n  22: [\"24\u00f78=\", \"89\u00f79=\"],\n  23: [\"41\u00f79=\", \"12\u00f75=\"],\n  24: [\"13\u00f73=\", \"71\u00f73=\"],\n  25: [\"52\u00f72=\", \"12\u00f76=\"],\n  41: [\"64\u00f73=\", \"59\u00f79=\"],\n  42: [\"41\u00f79=\", \"61\u00f78=\"],\n  43: [\"24\u00f78=\", \"28\u00f78=\"],\n  44: [\"76\u00f78=\", \"82\u00f76=\"],\n  45: [\"22\u00f77=\", \"36\u00f76=\"],\n  61: [\"73\u00f74=\", \"13\u00f72=\"],\n  62: [\"86\u00f73=\", \"74\u00f79=\"],\n  63: [\"64\u00f75=\", \"24\u00f74=\"],\n  64: [\"50\u00f78=\", \"74\u00f78=\"],\n  65: [\"49\u00f74=\", \"86\u00f72=\"],\n  81: [\"31\u00f79=\", \"85\u00f78=\"],\n  82: [\"70\u00f75=\", \"85\u00f74=\"],\n  83: [\"14\u00f78=\", \"88\u00f72=\"],\n  84: [\"85\u00f73=\", \"85\u00f77=\"],\n  85: [\"83\u00f76=\", \"18\u00f77=\"],\n};\n\nconst items = paragraphs.items;\nfor (const [idxStr, pair] of Object.entries(replacements)) {\n  const idx = Number(idxStr);\n  const [oldText, newText] = pair;\n  const para = items[idx];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Unexpected paragraph text at index \" + idx + \": \" + JSON.stringify(para.text)\n    );\n  }\n  // Replace the whole paragraph's text while keeping its (run)\n  // formatting intact.\n  para.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice\n# table. The table has 5 \"data\" rows (1, 5, 9, 13, 17) interleaved\n# with blank spacer rows, 5 columns each. We address cells by\n# row/column rather than searching for the text, because several of\n# the original problem strings repeat (e.g. \"86\u00f73=\", \"41\u00f79=\",\n# \"24\u00f78=\") and a simple Find/Replace-all would touch every occurrence\n# instead of the specific one the diff changes.\n\n$d = $word.ActiveDocument\n\n# --- Title date paragraph -------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -ne \"2024-01-23 Tuesday\") {\n    throw \"Unexpected title text: $($titlePara.Range.Text)\"\n}\n$titlePara.Range.Text = \"2024-01-24 Wednesday\"\n\n# --- Table of division problems -------------------------------------------\n$tbl = $d.Tables.Item(1)\n\n# Row number (1-based) -> ordered list of [old, new] pairs for the 5\n# cells in that row.\n$rowData = @{\n    1  = @(\n            @(\"93\u00f77=\", \"40\u00f79=\"),\n            @(\"86\u00f73=\", \"36\u00f77=\"),\n            @(\"43\u00f79=\", \"19\u00f75=\"),\n            @(\"98\u00f73=\", \"21\u00f78=\"),\n            @(\"78\u00f73=\", \"82\u00f72=\")\n         )\n    5  = @(\n            @(\"78\u00f78=\", \"74\u00f74=\"),\n            @(\"24\u00f78=\", \"89\u00f79=\"),\n            @(\"41\u00f79=\", \"12\u00f75=\"),\n            @(\"13\u00f73=\", \"71\u00f73=\"),\n            @(\"52\u00f72=\", \"12\u00f76=\")\n         )\n    9  = @(\n            @(\"64\u00f73=\", \"59\u00f79=\"),\n            @(\"41\u00f79=\", \"61\u00f78=\"),\n            @(\"24\u00f78=\", \"28\u00f78=\"),\n            @(\"76\u00f78=\", \"82\u00f76=\"),\n            @(\"22\u00f77=\", \"36\u00f76=\")\n         )\n    13 = @(\n            @(\"73\u00f74=\", \"13\u00f72=\"),\n            @(\"86\u00f73=\", \"74\u00f79=\"),\n            @(\"64\u00f75=\", \"24\u00f74=\"),\n            @(\"50\u00f78=\", \"74\u00f78=\"),\n            @(\"49\u00f74=\", \"86\u00f72=\")\n         )\n    17 = @(\n            @(\"31\u00f79=\", \"85\u00f78=\"),\n            @(\"70\u00f75=\", \"85\u00f74=\"),\n            @(\"14\u00f78=\", \"88\u00f72=\"),\n            @(\"85\u00f73=\", \"85\u00f77=\"),\n            @(\"83\u00f76=\", \"18\u00f77=\")\n         )\n}\n\nforeach ($rowIndex in $rowData.Keys) {\n    $row = $tbl.Rows.Item($rowIndex)\n    $cellPairs = $rowData[$rowIndex]\n    for ($c = 1; $c -le $cellPairs.Count; $c++) {\n        $cell = $row.Cells.Item($c)\n        $old = $cellPairs[$c - 1][0]\n        $new = $cellPairs[$c - 1][1]\n        $current = $cell.Range.Text.TrimEnd(\"`r\", \"`a\")\n        if ($current -ne $old) {\n            throw \"Unexpected text in row $rowIndex, col $c : $current\"\n        }\n        $cell.Range.Text = $new\n    }\n}\n"}
